$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (old E/F/G "Project/Fase/Code" shift right to F/G/H)
$ws.Columns("E:E").Insert()

# New header for the inserted column
$ws.Range("E1").Value2 = "MedewekerCode"

# Two new trailing columns
$ws.Range("I1").Value2 = "Medewerker"
$ws.Range("J1").Value2 = "Tijd-ID"

# Match column widths of the new columns as closely as possible
$ws.Columns("E:E").ColumnWidth = 15.67
$ws.Columns("I:I").ColumnWidth = 11.67

# Workbook-scoped defined names describing the data columns
$wb.Names.Add("IDs", "=Sheet1!`$A:`$A")
$wb.Names.Add("Dates", "=Sheet1!`$B:`$B")
$wb.Names.Add("TimeBegin", "=Sheet1!`$C:`$C")
$wb.Names.Add("TimeEnd", "=Sheet1!`$D:`$D")
$wb.Names.Add("EmployeeCodes", "=Sheet1!`$E:`$E")
$wb.Names.Add("Projects", "=Sheet1!`$F:`$F")
$wb.Names.Add("Phases", "=Sheet1!`$G:`$G")
$wb.Names.Add("Codes", "=Sheet1!`$H:`$H")
$wb.Names.Add("Employees", "=Sheet1!`$I:`$I")
$wb.Names.Add("TimeIDs", "=Sheet1!`$J:`$J")
$wb.Names.Add("ID", "=Sheet1!#REF!")

# Select the new column (as the user would after inserting it)
$null = $ws.Columns("E:E").Select()
